# Fill in the historical (actual) figures for the Base Case balance-sheet /
# cash-flow driver rows. These were previously blank input cells (B:F, the
# five historical years) feeding the projected columns (H:L) and the "Net"
# roll-forward column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base Case")

# Row -> historical values for columns B,C,D,E,F (an empty-string entry
# means "leave that historical cell blank", matching the source CSV which
# didn't have data that far back for some line items).
$rowData = [ordered]@{
    15 = @(11257, 11636, 11634, 11076, 10505)     # Capital Expenditures
    16 = @(-11488, -12229, -12352, -13168, -13548) # Depreciation & Amortization
    19 = @(791, 1764, 2111, 2349, 2132)            # Inventory
    20 = @(18383, 26287, 11233, 20481, 46671)      # Cash + Short Term Investments
    21 = @(57653, 73286, 68531, 89378, 106869)     # Current Assets
    22 = @(176064, 207000, 231839, 290345, 321686) # Total Assets
    23 = @($null, $null, 6308, 10999, 11605)       # Short Term Debt & Current Portion of LTD
    24 = @(38542, 43658, 63448, 80610, 79006)      # Current Liabilities
    25 = @($null, 16960, 28987, 53329, 75427)      # Total Long Term Debt
    26 = @(57854, 83451, 120292, 170990, 193437)   # Total Liabilities
    27 = @(118210, 123549, 111547, 119355, 128249) # Shareholders' Equity
}

$columns = @("B", "C", "D", "E", "F")

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $value = $values[$i]
        if ($null -ne $value) {
            $ws.Range("$($columns[$i])$row").Value = $value
        }
    }
}
